# Updated cryptos list on Fri Jul 26 20:49:29 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text representation
# (values like "579.25" or "67.535.86" must stay literal strings, not be
# reinterpreted as numbers/dates by Excel's automatic type inference).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '67.535.86'
$ws.Range("E2").Value = '  +3.99%  '
$ws.Range("D3").Value = '3.256.51'
$ws.Range("E3").Value = '  +4.05%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '579.25'
$ws.Range("E5").Value = '  +2.16%  '
$ws.Range("D6").Value = '181.48'
$ws.Range("E6").Value = '  +7.66%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("D9").Value = '3.255.37'
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("E10").Value = '  +8.78%  '
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("D12").Value = '0.415'
$ws.Range("E12").Value = '  +7.55%  '
$ws.Range("D13").Value = '3.822.44'
$ws.Range("E13").Value = '  +4.29%  '
$ws.Range("E14").Value = '  +1.59%  '
$ws.Range("E15").Value = '  +6.15%  '
$ws.Range("D16").Value = '67.504.58'
$ws.Range("E16").Value = '  +3.98%  '
$ws.Range("E17").Value = '  +4.59%  '
$ws.Range("D18").Value = '3.259.64'
$ws.Range("E18").Value = '  +4.43%  '
$ws.Range("E19").Value = '  +3.81%  '
$ws.Range("D20").Value = '13.51'
$ws.Range("E20").Value = '  +7.00%  '
$ws.Range("D21").Value = '375.06'
$ws.Range("E21").Value = '  +5.82%  '
$ws.Range("D22").Value = '7.62'
$ws.Range("E22").Value = '  +6.32%  '
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("D24").Value = '70.86'
$ws.Range("E24").Value = '  +3.54%  '
$ws.Range("E25").Value = '  +4.71%  '
$ws.Range("E26").Value = '  +7.41%  '
$ws.Range("D27").Value = '9.58'
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E28").Value = '  +3.79%  '
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  +4.25%  '
$ws.Range("E31").Value = '  +8.63%  '
$ws.Range("E32").Value = '  +5.33%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  +8.38%  '
$ws.Range("E35").Value = '  +6.46%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '1.51'
$ws.Range("E36").Value = '  +6.74%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '163.30'
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("E38").Value = '  +3.58%  '
$ws.Range("E39").Value = '  +6.32%  '
$ws.Range("D40").Value = '6.81'
$ws.Range("E40").Value = '  +13.20%  '
$ws.Range("D41").Value = '26.71'
$ws.Range("E41").Value = '  +3.16%  '
$ws.Range("E42").Value = '  +13.01%  '
$ws.Range("E43").Value = '  +7.83%  '
$ws.Range("D44").Value = '2.697.12'
$ws.Range("E44").Value = '  +2.96%  '
$ws.Range("D45").Value = '351.29'
$ws.Range("E45").Value = '  +9.93%  '
$ws.Range("D46").Value = '25.33'
$ws.Range("E46").Value = '  +7.99%  '
$ws.Range("D47").Value = '40.77'
$ws.Range("E47").Value = '  +3.51%  '
$ws.Range("E48").Value = '  +5.16%  '
$ws.Range("D49").Value = '0.0280'
$ws.Range("E49").Value = '  +3.90%  '
$ws.Range("D50").Value = '0.996'
$ws.Range("E50").Value = '  +7.36%  '
$ws.Range("E51").Value = '  +1.43%  '
